# New Course Outline is being prepared
# Clear the roster of student names/IDs from Sheet1 and Sheet2, leaving
# the lookup formulas in place (they will resolve to #N/A / 0 via IFERROR).

$wb = $excel.ActiveWorkbook

# --- Sheet1: clear column A (Name) for every populated data row ---------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A2:A5").ClearContents()
$ws1.Range("A7:A10").ClearContents()
$ws1.Range("A12:A15").ClearContents()
$ws1.Range("A17:A20").ClearContents()
$ws1.Range("A22:A25").ClearContents()
$ws1.Range("A27:A30").ClearContents()

# --- Sheet2: clear columns A (Name) and B (ID) for every populated row --
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A2:B25").ClearContents()

# --- View/selection bookkeeping ------------------------------------------
# Sheet2 keeps a stored selection of A2:B25 but is no longer the active tab.
$ws2.Range("A2:B25").Select() | Out-Null

# Sheet1 becomes the active tab with A2 selected.
$ws1.Activate()
$ws1.Range("A2").Select() | Out-Null
$win = $excel.ActiveWindow
$win.Zoom = 220
